$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item("Subtitle 3")
$sh.TextFrame.TextRange.Text = "infoway"
